# Update "想去人数" (attendee interest count) figures in column F across
# the refreshed convention-listing sheets (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# 展览 (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$map1 = @{
    2  = 549
    3  = 247
    6  = 673
    7  = 358
    8  = 39
    10 = 428
    11 = 6334
    12 = 117
    13 = 26
    14 = 1895
    15 = 4739
    19 = 5520
    20 = 7303
    22 = 1092
    23 = 763
    24 = 4037
    25 = 565
    26 = 74
    29 = 1069
    31 = 567
    32 = 696
    33 = 1703
    34 = 242
    35 = 1918
    36 = 236
    37 = 45
    38 = 1252
    39 = 1334
    40 = 697
    42 = 1349
    43 = 3708
    45 = 348
    48 = 97
}
foreach ($row in $map1.Keys) {
    $ws1.Range("F$row").Value = $map1[$row]
}

# 演出 (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$map2 = @{
    3 = 1277
}
foreach ($row in $map2.Keys) {
    $ws2.Range("F$row").Value = $map2[$row]
}

# 本地生活 (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$map3 = @{
    2 = 4465
}
foreach ($row in $map3.Keys) {
    $ws3.Range("F$row").Value = $map3[$row]
}

# 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$map4 = @{
    2  = 4465
    3  = 549
    4  = 1277
    6  = 247
    11 = 673
    12 = 358
    14 = 428
    15 = 117
    16 = 26
    17 = 4739
    18 = 5520
    19 = 5520
    21 = 1092
    22 = 763
    23 = 4037
    24 = 565
    28 = 1069
    30 = 567
    31 = 696
    32 = 1703
    33 = 242
    34 = 1918
    39 = 697
    43 = 3708
    46 = 348
    48 = 97
}
foreach ($row in $map4.Keys) {
    $ws4.Range("F$row").Value = $map4[$row]
}
